$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 22.700661
$ws.Range("H2").Value = 68.10198299999999
$ws.Range("I2").Value = 0.08615268874617349
$ws.Range("J2").Value = 0.08615268874617349
$ws.Range("M2").Value = 2.565830333333333
$ws.Range("N2").Value = 7.697490999999999
$ws.Range("O2").Value = 0.0934185609347503
$ws.Range("P2").Value = 0.0934185609347503
$ws.Range("Q2").Value = 58.24604458051699
$ws.Range("R2").Value = 524.2144012246529
$ws.Range("S2").Value = 0.008048260203326985
$ws.Range("T2").Value = 0.008048260203326985
$ws.Range("G3").Value = 22.700661
$ws.Range("H3").Value = 68.10198299999999
$ws.Range("I3").Value = 0.08615268874617349
$ws.Range("J3").Value = 0.08615268874617349
$ws.Range("O3").Value = 0.3847798091300315
$ws.Range("P3").Value = 0.3847798091300315
$ws.Range("Q3").Value = 239.9084474435929
$ws.Range("R3").Value = 2159.176026992337
$ws.Range("S3").Value = 0.03314981513179165
$ws.Range("T3").Value = 0.03314981513179165
$ws.Range("G4").Value = 22.700661
$ws.Range("H4").Value = 68.10198299999999
$ws.Range("I4").Value = 0.08615268874617349
$ws.Range("J4").Value = 0.08615268874617349
$ws.Range("M4").Value = 13.68376133333333
$ws.Range("N4").Value = 41.051284
$ws.Range("O4").Value = 0.4982080363333638
$ws.Range("P4").Value = 0.4982080363333638
$ws.Range("Q4").Value = 310.630427232908
$ws.Range("R4").Value = 2795.673845096172
$ws.Range("S4").Value = 0.04292196188507059
$ws.Range("T4").Value = 0.04292196188507059
$ws.Range("G5").Value = 22.700661
$ws.Range("H5").Value = 68.10198299999999
$ws.Range("I5").Value = 0.08615268874617349
$ws.Range("J5").Value = 0.08615268874617349
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.6480206666666667
$ws.Range("N5").Value = 1.944062
$ws.Range("O5").Value = 0.0235935936018545
$ws.Range("P5").Value = 0.0235935936018545
$ws.Range("Q5").Value = 14.710497474994
$ws.Range("R5").Value = 132.394477274946
$ws.Range("S5").Value = 0.002032651525984281
$ws.Range("T5").Value = 0.002032651525984281
$ws.Range("I6").Value = 0.5030288587986086
$ws.Range("J6").Value = 0.5030288587986087
$ws.Range("M6").Value = 2.565830333333333
$ws.Range("N6").Value = 7.697490999999999
$ws.Range("O6").Value = 0.0934185609347503
$ws.Range("P6").Value = 0.0934185609347503
$ws.Range("Q6").Value = 340.0873700087705
$ws.Range("R6").Value = 3060.786330078934
$ws.Range("S6").Value = 0.04699223209761572
$ws.Range("T6").Value = 0.04699223209761573
$ws.Range("I7").Value = 0.5030288587986086
$ws.Range("J7").Value = 0.5030288587986087
$ws.Range("O7").Value = 0.3847798091300315
$ws.Range("P7").Value = 0.3847798091300315
$ws.Range("S7").Value = 0.1935553482754262
$ws.Range("T7").Value = 0.1935553482754262
$ws.Range("I8").Value = 0.5030288587986086
$ws.Range("J8").Value = 0.5030288587986087
$ws.Range("M8").Value = 13.68376133333333
$ws.Range("N8").Value = 41.051284
$ws.Range("O8").Value = 0.4982080363333638
$ws.Range("P8").Value = 0.4982080363333638
$ws.Range("Q8").Value = 1813.710884630216
$ws.Range("R8").Value = 16323.39796167194
$ws.Range("S8").Value = 0.2506130199610677
$ws.Range("T8").Value = 0.2506130199610678
$ws.Range("I9").Value = 0.5030288587986086
$ws.Range("J9").Value = 0.5030288587986087
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.6480206666666667
$ws.Range("N9").Value = 1.944062
$ws.Range("O9").Value = 0.0235935936018545
$ws.Range("P9").Value = 0.0235935936018545
$ws.Range("Q9").Value = 85.89174481840777
$ws.Range("R9").Value = 773.02570336567
$ws.Range("S9").Value = 0.01186825846449902
$ws.Range("T9").Value = 0.01186825846449902
$ws.Range("G10").Value = 41.94534433333333
$ws.Range("H10").Value = 125.836033
$ws.Range("I10").Value = 0.159189381961201
$ws.Range("J10").Value = 0.159189381961201
$ws.Range("M10").Value = 2.565830333333333
$ws.Range("N10").Value = 7.697490999999999
$ws.Range("O10").Value = 0.0934185609347503
$ws.Range("P10").Value = 0.0934185609347503
$ws.Range("Q10").Value = 107.6246368325781
$ws.Range("R10").Value = 968.6217314932028
$ws.Range("S10").Value = 0.0148712429789077
$ws.Range("T10").Value = 0.0148712429789077
$ws.Range("G11").Value = 41.94534433333333
$ws.Range("H11").Value = 125.836033
$ws.Range("I11").Value = 0.159189381961201
$ws.Range("J11").Value = 0.159189381961201
$ws.Range("O11").Value = 0.3847798091300315
$ws.Range("P11").Value = 0.3847798091300315
$ws.Range("Q11").Value = 443.2929259855874
$ws.Range("R11").Value = 3989.636333870286
$ws.Range("S11").Value = 0.06125286000655859
$ws.Range("T11").Value = 0.06125286000655859
$ws.Range("G12").Value = 41.94534433333333
$ws.Range("H12").Value = 125.836033
$ws.Range("I12").Value = 0.159189381961201
$ws.Range("J12").Value = 0.159189381961201
$ws.Range("M12").Value = 13.68376133333333
$ws.Range("N12").Value = 41.051284
$ws.Range("O12").Value = 0.4982080363333638
$ws.Range("P12").Value = 0.4982080363333638
$ws.Range("Q12").Value = 573.9700809018192
$ws.Range("R12").Value = 5165.730728116372
$ws.Range("S12").Value = 0.07930942939201176
$ws.Range("T12").Value = 0.07930942939201176
$ws.Range("G13").Value = 41.94534433333333
$ws.Range("H13").Value = 125.836033
$ws.Range("I13").Value = 0.159189381961201
$ws.Range("J13").Value = 0.159189381961201
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.6480206666666667
$ws.Range("N13").Value = 1.944062
$ws.Range("O13").Value = 0.0235935936018545
$ws.Range("P13").Value = 0.0235935936018545
$ws.Range("Q13").Value = 27.18144999844955
$ws.Range("R13").Value = 244.633049986046
$ws.Range("S13").Value = 0.003755849583722964
$ws.Range("T13").Value = 0.003755849583722964
$ws.Range("G14").Value = 66.302588
$ws.Range("H14").Value = 198.907764
$ws.Range("I14").Value = 0.2516290704940168
$ws.Range("J14").Value = 0.2516290704940168
$ws.Range("M14").Value = 2.565830333333333
$ws.Range("N14").Value = 7.697490999999999
$ws.Range("O14").Value = 0.0934185609347503
$ws.Range("P14").Value = 0.0934185609347503
$ws.Range("Q14").Value = 170.1211914689027
$ws.Range("R14").Value = 1531.090723220124
$ws.Range("S14").Value = 0.02350682565489989
$ws.Range("T14").Value = 0.02350682565489989
$ws.Range("G15").Value = 66.302588
$ws.Range("H15").Value = 198.907764
$ws.Range("I15").Value = 0.2516290704940168
$ws.Range("J15").Value = 0.2516290704940168
$ws.Range("O15").Value = 0.3847798091300315
$ws.Range("P15").Value = 0.3847798091300315
$ws.Range("Q15").Value = 700.7087127803106
$ws.Range("R15").Value = 6306.378415022796
$ws.Range("S15").Value = 0.09682178571625503
$ws.Range("T15").Value = 0.09682178571625503
$ws.Range("G16").Value = 66.302588
$ws.Range("H16").Value = 198.907764
$ws.Range("I16").Value = 0.2516290704940168
$ws.Range("J16").Value = 0.2516290704940168
$ws.Range("M16").Value = 13.68376133333333
$ws.Range("N16").Value = 41.051284
$ws.Range("O16").Value = 0.4982080363333638
$ws.Range("P16").Value = 0.4982080363333638
$ws.Range("Q16").Value = 907.2687899743307
$ws.Range("R16").Value = 8165.419109768976
$ws.Range("S16").Value = 0.1253636250952137
$ws.Range("T16").Value = 0.1253636250952137
$ws.Range("G17").Value = 66.302588
$ws.Range("H17").Value = 198.907764
$ws.Range("I17").Value = 0.2516290704940168
$ws.Range("J17").Value = 0.2516290704940168
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.6480206666666667
$ws.Range("N17").Value = 1.944062
$ws.Range("O17").Value = 0.0235935936018545
$ws.Range("P17").Value = 0.0235935936018545
$ws.Range("Q17").Value = 42.96544727748533
$ws.Range("R17").Value = 386.6890254973679
$ws.Range("S17").Value = 0.005936834027648229
$ws.Range("T17").Value = 0.005936834027648229

Write-Host "Applied 182 cell updates"
